$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (Result / D1) ---------------------------------------
$ws.Range("D1").Value = "Result"
$ws.Range("E1").Value = "D1"

# --- Row 2: 1 + 7 = 8 -------------------------------------------------------
# Leading apostrophe forces text storage (so "8", "10", "27"... don't get
# reinterpreted as numbers) without disturbing the General number format.
$ws.Range("D2").Value = "'8`t`n"
$ws.Range("E2").Value = "D2"

# --- Row 3: 2 + 8 = 10 -------------------------------------------------------
$ws.Range("D3").Value = "'10`t`n"
$ws.Range("E3").Value = "D3"

# --- Row 4: 3 * 9 = 27 -------------------------------------------------------
$ws.Range("D4").Value = "'27`t`n"

# --- Row 5: 4 * 0 = 0 -------------------------------------------------------
$ws.Range("D5").Value = "'0`t`n"

# --- Row 6: 5 * 1 = 5 -------------------------------------------------------
$ws.Range("D6").Value = "'5`t`n"

# --- Row 7: 6 - 2 = 4 -------------------------------------------------------
$ws.Range("D7").Value = "'4`t`n"

# Entering text with an embedded newline auto-grows a row's height; put rows
# 3-7 back to the sheet default before fixing row 2 at its real height.
$ws.Range("D3:D7").EntireRow.AutoFit()

# --- Formatting: D2's two-line text wraps and row 2 grows to fit it --------
$ws.Range("D2").WrapText = $true
$ws.Rows.Item(2).RowHeight = 30

# --- Column D width (auto best-fit like Excel would compute) ---------------
$ws.Columns.Item(4).AutoFit()

# --- Selection matches the saved cursor position ----------------------------
$ws.Range("D7").Select()
